# The "LOB1207 -  Poluição Ambiental I  (Requisito fraco)" requisito row
# (row 25, under the "Requisitos:" section) was removed from the sheet.
# Deleting the entire row shifts the following row ("LOQ4233 -  Gestão de
# Negócios  (Requisito fraco)") up into its place, shrinks the used range
# by one row, and drops the now-unreferenced shared string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(25).Delete()
